$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column E. This shifts the former columns
# E:Q (headers ESTADO..BANCO, the PAGOS/MOTIVO blocks, etc.) three slots
# to the right, becoming H:T, and adjusts merged cell ranges / dimension
# automatically.
$ws.Range("E1:G1").EntireColumn.Insert()

# Give the three freshly inserted columns (E:G) their new width.
$ws.Range("E1:G1").ColumnWidth = 23

# Populate the header row (row 10) for the new columns.
$ws.Range("E10").Value = "TIP COMPROBANTE"
$ws.Range("F10").Value = "DOC ASOCIADO"
$ws.Range("G10").Value = "FECHA EMISIÓN"

# Match the saved selection/active cell.
$ws.Range("G10").Select() | Out-Null
